$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.09982999999999999
$ws.Range("H2").Value = 0.29949
$ws.Range("I2").Value = 0.06953924013247029
$ws.Range("J2").Value = 0.06953924013247029
$ws.Range("M2").Value = 6.045145666666667
$ws.Range("N2").Value = 18.135437
$ws.Range("O2").Value = 0.8160840232643366
$ws.Range("P2").Value = 0.8160840232643367
$ws.Range("Q2").Value = 0.6034868919033334
$ws.Range("R2").Value = 5.431382027130001
$ws.Range("S2").Value = 0.05674986286205117
$ws.Range("T2").Value = 0.05674986286205118

$ws.Range("G3").Value = 0.09982999999999999
$ws.Range("H3").Value = 0.29949
$ws.Range("I3").Value = 0.06953924013247029
$ws.Range("J3").Value = 0.06953924013247029
$ws.Range("O3").Value = 0.09212864864242169
$ws.Range("P3").Value = 0.09212864864242169
$ws.Range("Q3").Value = 0.06812831796666667
$ws.Range("R3").Value = 0.6131548617
$ws.Range("S3").Value = 0.006406556221025345
$ws.Range("T3").Value = 0.006406556221025345

$ws.Range("G4").Value = 0.09982999999999999
$ws.Range("H4").Value = 0.29949
$ws.Range("I4").Value = 0.06953924013247029
$ws.Range("J4").Value = 0.06953924013247029
$ws.Range("M4").Value = 0.6799149999999999
$ws.Range("N4").Value = 2.039745
$ws.Range("O4").Value = 0.09178732809324164
$ws.Range("P4").Value = 0.09178732809324165
$ws.Range("Q4").Value = 0.06787591444999999
$ws.Range("R4").Value = 0.6108832300499999
$ws.Range("S4").Value = 0.006382821049393766
$ws.Range("T4").Value = 0.006382821049393767

$ws.Range("I5").Value = 0.4393303855760352
$ws.Range("J5").Value = 0.4393303855760352
$ws.Range("M5").Value = 6.045145666666667
$ws.Range("N5").Value = 18.135437
$ws.Range("O5").Value = 0.8160840232643366
$ws.Range("P5").Value = 0.8160840232643367
$ws.Range("Q5").Value = 3.812669341869556
$ws.Range("R5").Value = 34.31402407682601
$ws.Range("S5").Value = 0.3585305086031631
$ws.Range("T5").Value = 0.3585305086031632

$ws.Range("I6").Value = 0.4393303855760352
$ws.Range("J6").Value = 0.4393303855760352
$ws.Range("O6").Value = 0.09212864864242169
$ws.Range("P6").Value = 0.09212864864242169
$ws.Range("S6").Value = 0.0404749147306742
$ws.Range("T6").Value = 0.0404749147306742

$ws.Range("I7").Value = 0.4393303855760352
$ws.Range("J7").Value = 0.4393303855760352
$ws.Range("M7").Value = 0.6799149999999999
$ws.Range("N7").Value = 2.039745
$ws.Range("O7").Value = 0.09178732809324164
$ws.Range("P7").Value = 0.09178732809324165
$ws.Range("Q7").Value = 0.4288219372233333
$ws.Range("R7").Value = 3.85939743501
$ws.Range("S7").Value = 0.0403249622421979
$ws.Range("T7").Value = 0.04032496224219791

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.705063
$ws.Range("H8").Value = 2.115189
$ws.Range("I8").Value = 0.4911303742914945
$ws.Range("J8").Value = 0.4911303742914945
$ws.Range("M8").Value = 6.045145666666667
$ws.Range("N8").Value = 18.135437
$ws.Range("O8").Value = 0.8160840232643366
$ws.Range("P8").Value = 0.8160840232643367
$ws.Range("Q8").Value = 4.262208539177
$ws.Range("R8").Value = 38.359876852593
$ws.Range("S8").Value = 0.4008036517991224
$ws.Range("T8").Value = 0.4008036517991224

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.705063
$ws.Range("H9").Value = 2.115189
$ws.Range("I9").Value = 0.4911303742914945
$ws.Range("J9").Value = 0.4911303742914945
$ws.Range("O9").Value = 0.09212864864242169
$ws.Range("P9").Value = 0.09212864864242169
$ws.Range("Q9").Value = 0.4811655439300001
$ws.Range("R9").Value = 4.33048989537
$ws.Range("S9").Value = 0.04524717769072215
$ws.Range("T9").Value = 0.04524717769072215

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.705063
$ws.Range("H10").Value = 2.115189
$ws.Range("I10").Value = 0.4911303742914945
$ws.Range("J10").Value = 0.4911303742914945
$ws.Range("M10").Value = 0.6799149999999999
$ws.Range("N10").Value = 2.039745
$ws.Range("O10").Value = 0.09178732809324164
$ws.Range("P10").Value = 0.09178732809324165
$ws.Range("Q10").Value = 0.479382909645
$ws.Range("R10").Value = 4.314446186805
$ws.Range("S10").Value = 0.04507954480164997
$ws.Range("T10").Value = 0.04507954480164998
